$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.717.05"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "3.312.04"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.310.47"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "3.858.39"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.312.22"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "63.789.75"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +6.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("E30").Value = "  -5.49%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.39%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "0.0₃0747"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "434.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("D41").Value = "3.114.76"
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("E42").Value = "  +9.06%  "
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.12%  "

Write-Host "Applied 71 cell updates."
